# Applies the "Atualização automática" data refresh to Sheet1:
# - New random UUIDs in column A for rows 2-19
# - Updated detection image filename / bounding box / confidence values
#   for a handful of rows (15 and 16) that were re-processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A (Fly_ID) - new UUIDs for rows 2-19
$ws.Range("A2").Value  = "320891a7-9716-4af5-a792-4bc855814887"
$ws.Range("A3").Value  = "cf5b1d2d-e6a6-4890-86b7-727104c438c7"
$ws.Range("A4").Value  = "d3e26aa8-f8ad-4659-9e10-e2504a577d64"
$ws.Range("A5").Value  = "4b7d1a0a-7be9-4095-8d08-3fe1aa818667"
$ws.Range("A6").Value  = "2cdb1266-8286-4b3d-9247-3002269e62d9"
$ws.Range("A7").Value  = "64c4f383-9171-430f-bdb5-9a39599a9a14"
$ws.Range("A8").Value  = "1cccb561-3eb4-4319-b0e2-f02fff0cc1f5"
$ws.Range("A9").Value  = "78fef29f-aac1-48da-8b16-c80cb5477c44"
$ws.Range("A10").Value = "30471ab7-ea71-480e-a311-459c261009c6"
$ws.Range("A11").Value = "d69af814-f329-4aa8-af85-42f1510c83c3"
$ws.Range("A12").Value = "66d9d4d1-2522-46ae-8e8c-74cf83b0c8a9"
$ws.Range("A13").Value = "0cf2d1a7-05a1-47da-b6fd-1d3d8e821560"
$ws.Range("A14").Value = "5cf5c49a-8520-4dc8-b6a5-615b23a4cd37"
$ws.Range("A15").Value = "bc03209c-34a6-4ef6-8e16-c8dd11e91222"
$ws.Range("A16").Value = "6e29e66a-19b5-435f-86f1-621f4e1e2d90"
$ws.Range("A17").Value = "c86d904d-e99d-4a2c-b3b8-0a7633148ae1"
$ws.Range("A18").Value = "cce28004-92b0-486c-8e20-90e6607453bc"
$ws.Range("A19").Value = "b72ae262-fded-42d8-8835-c03ab20883a9"

# Row 15 - detection image refined, bounding box tweaked
$ws.Range("D15").Value = "image_20250807111344_ppp0.jpg"

# I15's new value ("794,481,831,526") reads like a number with thousands
# separators, so Excel would auto-convert it to a numeric cell. Force the
# cell to stay text (matching the original inline-string coordinate list)
# by pre-formatting as Text, then restoring the default style afterward.
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = "794,481,831,526"
$ws.Range("I15").Style = "Normal"

# Row 16 - detection image refined, bounding box + confidence tweaked
$ws.Range("D16").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I16").Value = "1182,409,1232,451"

# J16's new value ("0.75") reads like a plain decimal number, so guard it
# the same way to keep it stored as text rather than a numeric value.
$ws.Range("J16").NumberFormat = "@"
$ws.Range("J16").Value = "0.75"
$ws.Range("J16").Style = "Normal"
